# Regenerate orders with updated distance/sizes.
# Apply text substitutions across the used range of the active sheet:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# These substrings appear embedded inside larger strings (Condition,
# Filename_Left, Filename_Right, Distance, Size columns), so we do a
# substring replace on every string cell in the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

$replacements = @(
    @("D64", "D69"),
    @("D80", "D86"),
    @("D51", "D55"),
    @("S30", "S31")
)

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()

        if ($val -is [string]) {
            $newVal = $val
            foreach ($pair in $replacements) {
                $newVal = $newVal.Replace($pair[0], $pair[1])
            }
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
